$wb = $excel.ActiveWorkbook

$sheetEnglish = $wb.Worksheets.Item("english")

# Set D4 on the "english" sheet to a new "Video off" string (adds a new shared
# string entry and points D4 at it, replacing the previous "Video Off" value).
$sheetEnglish.Range("D4").Value = "Video off"

# Move the selection on the "english" sheet to D4.
$sheetEnglish.Range("D4").Select()

# Make the "english" sheet the active tab (this also updates the workbook's
# active-tab/selected-sheet bookkeeping, moving tabSelected off "hindi").
$sheetEnglish.Activate()
